# Remove "<strong>" / "</strong>" wrapper tags from the "[sinopsis]" header
# labels in the first worksheet. The trailing space before the (removed)
# closing tag is preserved, matching the target shared-string values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$replacements = @{
    "A2"  = "Activos [sinopsis] "
    "A3"  = "Activos corrientes [sinopsis] "
    "A15" = "Activos no corrientes [sinopsis] "
    "A31" = "Patrimonio y pasivos [sinopsis] "
    "A32" = "Pasivos [sinopsis] "
    "A33" = "Pasivos corrientes [sinopsis] "
    "A44" = "Pasivos no corrientes [sinopsis] "
    "A55" = "Patrimonio [sinopsis] "
}

foreach ($addr in $replacements.Keys) {
    $ws.Range($addr).Value = $replacements[$addr]
}
